{"js": "// Apply the certificate-template edits:\n//   1) serial number: 1015711121811 -> 1015705011901\n//   2) \"valid until\" date: 10 \u0433\u0440\u0443\u0434\u043d\u044f 2022 \u0440. -> 5 c\u0456\u0447\u043d\u044f 2023 \u0440.\n//   3) task/order number placeholder: \" \u043d\u043e\u043c\u0435\u0440 \u0437\u0430\u0432\u0434\u0430\u043d\u043d\u044f\" -> \"122627\"\n//   4) stamp date: 10 \u0433\u0440\u0443\u0434\u043d\u044f 2018 \u0440. -> 5 c\u0456\u0447\u043d\u044f 2019 \u0440.\n\nconst body = context.document.body;\n\nconst replacements = [\n  [\"1015711121811\", \"1015705011901\"],\n  [\"10 \u0433\u0440\u0443\u0434\u043d\u044f 2022 \u0440.\", \"5 c\u0456\u0447\u043d\u044f 2023 \u0440.\"],\n  [\" \u043d\u043e\u043c\u0435\u0440 \u0437\u0430\u0432\u0434\u0430\u043d\u043d\u044f\", \"122627\"],\n  [\"10 \u0433\u0440\u0443\u0434\u043d\u044f 2018 \u0440.\", \"5 c\u0456\u0447\u043d\u044f 2019 \u0440.\"]\n];\n\nconst allResults = replacements.map(([find]) =>\n  body.search(find, { matchCase: true, matchWholeWord: false })\n);\nallResults.forEach(r => r.load(\"items,text\"));\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, replaceWith] = replacements[i];\n  const results = allResults[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(replaceWith, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the certificate-template edits:\n#   1) serial number: 1015711121811 -> 1015705011901\n#   2) \"valid until\" date: 10 \u0433\u0440\u0443\u0434\u043d\u044f 2022 \u0440. -> 5 c\u0456\u0447\u043d\u044f 2023 \u0440.\n#   3) task/order number placeholder: \" \u043d\u043e\u043c\u0435\u0440 \u0437\u0430\u0432\u0434\u0430\u043d\u043d\u044f\" -> \"122627\"\n#   4) stamp date: 10 \u0433\u0440\u0443\u0434\u043d\u044f 2018 \u0440. -> 5 c\u0456\u0447\u043d\u044f 2019 \u0440.\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\nReplace-Text \"1015711121811\" \"1015705011901\"\nReplace-Text \"10 \u0433\u0440\u0443\u0434\u043d\u044f 2022 \u0440.\" \"5 c\u0456\u0447\u043d\u044f 2023 \u0440.\"\nReplace-Text \" \u043d\u043e\u043c\u0435\u0440 \u0437\u0430\u0432\u0434\u0430\u043d\u043d\u044f\" \"122627\"\nReplace-Text \"10 \u0433\u0440\u0443\u0434\u043d\u044f 2018 \u0440.\" \"5 c\u0456\u0447\u043d\u044f 2019 \u0440.\"\n"}
